$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '61.895.61'
Set-TextCell $ws.Range("E2") '  +0.55%  '

# Row 3
Set-TextCell $ws.Range("D3") '3.413.13'
Set-TextCell $ws.Range("E3") '  +0.78%  '

# Row 4
Set-TextCell $ws.Range("E4") '  +0.18%  '

# Row 5
Set-TextCell $ws.Range("D5") '409.38'
Set-TextCell $ws.Range("E5") '  +0.33%  '

# Row 6
Set-TextCell $ws.Range("D6") '128.14'
Set-TextCell $ws.Range("E6") '  -5.41%  '

# Row 7
Set-TextCell $ws.Range("D7") '0.621'
Set-TextCell $ws.Range("E7") '  +4.23%  '

# Row 8
Set-TextCell $ws.Range("D8") '1.00'
Set-TextCell $ws.Range("E8") '  -0.04%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.748'
Set-TextCell $ws.Range("E9") '  +10.72%  '

# Row 10
Set-TextCell $ws.Range("D10") '0.140'
Set-TextCell $ws.Range("E10") '  +15.66%  '

# Row 11
Set-TextCell $ws.Range("D11") '42.81'
Set-TextCell $ws.Range("E11") '  -0.01%  '

# Row 12
Set-TextCell $ws.Range("D12") '0.140'
Set-TextCell $ws.Range("E12") '  -0.59%  '

# Row 13
Set-TextCell $ws.Range("B13") 'Chainlink'
Set-TextCell $ws.Range("C13") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws.Range("D13") '20.96'
Set-TextCell $ws.Range("E13") '  +5.95%  '

# Row 14
Set-TextCell $ws.Range("B14") 'Polkadot'
Set-TextCell $ws.Range("C14") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range("D14") '8.81'
Set-TextCell $ws.Range("E14") '  +4.42%  '

# Row 15
Set-TextCell $ws.Range("D15") '0.0000198'
Set-TextCell $ws.Range("E15") '  +54.79%  '

# Row 16
Set-TextCell $ws.Range("D16") '3.417.04'
Set-TextCell $ws.Range("E16") '  +0.79%  '

# Row 17
Set-TextCell $ws.Range("D17") '12.66'
Set-TextCell $ws.Range("E17") '  +14.52%  '

# Row 18
Set-TextCell $ws.Range("D18") '1.06'
Set-TextCell $ws.Range("E18") '  +3.75%  '

# Row 19
Set-TextCell $ws.Range("D19") '61.981.58'
Set-TextCell $ws.Range("E19") '  +0.77%  '

# Row 20
Set-TextCell $ws.Range("D20") '401.85'
Set-TextCell $ws.Range("E20") '  +27.21%  '

# Row 21
Set-TextCell $ws.Range("D21") '89.97'
Set-TextCell $ws.Range("E21") '  +5.56%  '

# Row 22
Set-TextCell $ws.Range("D22") '3.18'
Set-TextCell $ws.Range("E22") '  -1.20%  '

# Row 23
Set-TextCell $ws.Range("D23") '13.35'
Set-TextCell $ws.Range("E23") '  +4.09%  '

# Row 24
Set-TextCell $ws.Range("D24") '3.24'
Set-TextCell $ws.Range("E24") '  +2.84%  '

# Row 25
Set-TextCell $ws.Range("D25") '32.97'
Set-TextCell $ws.Range("E25") '  +11.31%  '

# Row 26
Set-TextCell $ws.Range("B26") 'LEO'
Set-TextCell $ws.Range("C26") 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws.Range("D26") '4.79'
Set-TextCell $ws.Range("E26") '  +0.17%  '

# Row 27
Set-TextCell $ws.Range("B27") 'Filecoin'
Set-TextCell $ws.Range("C27") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D27") '8.50'
Set-TextCell $ws.Range("E27") '  +1.61%  '

# Row 28
Set-TextCell $ws.Range("D28") '7.63'
Set-TextCell $ws.Range("E28") '  -0.53%  '

# Row 29
Set-TextCell $ws.Range("D29") '2.74'
Set-TextCell $ws.Range("E29") '  +7.41%  '

# Row 30
Set-TextCell $ws.Range("D30") '0.117'
Set-TextCell $ws.Range("E30") '  -0.50%  '

# Row 31
Set-TextCell $ws.Range("B31") 'Kaspa'
Set-TextCell $ws.Range("C31") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range("D31") '0.171'
Set-TextCell $ws.Range("E31") '  -0.82%  '

# Row 32
Set-TextCell $ws.Range("B32") 'InjectiveProtocol'
Set-TextCell $ws.Range("C32") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range("D32") '43.72'
Set-TextCell $ws.Range("E32") '  +7.01%  '

# Row 33
Set-TextCell $ws.Range("D33") '11.77'
Set-TextCell $ws.Range("E33") '  +3.46%  '

# Row 34
Set-TextCell $ws.Range("E34") '  -0.06%  '

# Row 35
Set-TextCell $ws.Range("D35") '0.0498'
Set-TextCell $ws.Range("E35") '  +3.20%  '

# Row 36
Set-TextCell $ws.Range("D36") '52.51'
Set-TextCell $ws.Range("E36") '  +1.11%  '

# Row 37
Set-TextCell $ws.Range("E37") '  +0.09%  '

# Row 38
Set-TextCell $ws.Range("D38") '3.37'
Set-TextCell $ws.Range("E38") '  -1.95%  '

# Row 39
Set-TextCell $ws.Range("B39") 'Stacks'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws.Range("D39") '2.89'
Set-TextCell $ws.Range("E39") '  -1.51%  '

# Row 40
Set-TextCell $ws.Range("B40") 'Stellar'
Set-TextCell $ws.Range("C40") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range("D40") '0.131'
Set-TextCell $ws.Range("E40") '  +5.11%  '

# Row 41
Set-TextCell $ws.Range("D41") '0.313'
Set-TextCell $ws.Range("E41") '  +6.35%  '

# Row 42
Set-TextCell $ws.Range("D42") '140.33'
Set-TextCell $ws.Range("E42") '  +1.23%  '

# Row 43
Set-TextCell $ws.Range("D43") '1.97'
Set-TextCell $ws.Range("E43") '  -0.67%  '

# Row 44
Set-TextCell $ws.Range("D44") '4.01'
Set-TextCell $ws.Range("E44") '  -0.52%  '

# Row 45
Set-TextCell $ws.Range("D45") '2.38'
Set-TextCell $ws.Range("E45") '  +7.01%  '

# Row 46
Set-TextCell $ws.Range("D46") '16.69'
Set-TextCell $ws.Range("E46") '  -0.55%  '

# Row 47
Set-TextCell $ws.Range("D47") '21.87'
Set-TextCell $ws.Range("E47") '  +1.66%  '

# Row 48
Set-TextCell $ws.Range("D48") '2.102.29'
Set-TextCell $ws.Range("E48") '  -1.39%  '

# Row 49
Set-TextCell $ws.Range("E49") '  +0.26%  '

# Row 50
Set-TextCell $ws.Range("B50") 'BEAM'
Set-TextCell $ws.Range("C50") 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
Set-TextCell $ws.Range("D50") '0.0371'
Set-TextCell $ws.Range("E50") '  +8.07%  '

# Row 51
Set-TextCell $ws.Range("B51") 'Cronos'
Set-TextCell $ws.Range("C51") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range("D51") '0.126'
Set-TextCell $ws.Range("E51") '  +12.88%  '
